# Update maternal side of EMPDSR
# Rename the icd_mmg classification codes in column F (and keep icd_mmg_desc
# in column J unchanged) from the old "G1".."G9" scheme to the new
# "MM1".."MM9" scheme for the maternal mortality group rows (918-1129).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 918; $r -le 1129; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F = icd_mmg
    $val = $cell.Value2
    if ($val -match "^G(\d+)$") {
        $num = $matches[1]
        $cell.Value = "MM" + $num
    }
}

# Reflect the cursor/viewport position the author ended up at after the edit
$ws.Range("C1129").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1103
